$wb = $excel.ActiveWorkbook

$wsSpareReq = $wb.Worksheets.Item("spare_parts_required")
$wsSpareReq.Range("B2").Value = 3
$wsSpareReq.Range("C2").Value = 2
$wsSpareReq.Range("C3").Value = 4

$wsVessels = $wb.Worksheets.Item("vessels")
$wsVessels.Range("L2").Value = 36
$wsVessels.Range("J5").Value = 10000
$wsVessels.Select()
$wsVessels.Range("M9").Select()

$wsSpareReq.Select()
$wsSpareReq.Range("F5").Select()
